$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Fecha" column (D) values — shift weekly dates as described in the diff.
$ws.Range("D2").Value2 = 44832
$ws.Range("D3").Value2 = 44832
$ws.Range("D4").Value2 = 44846
$ws.Range("D5").Value2 = 44846
$ws.Range("D6").Value2 = 44838
$ws.Range("D7").Value2 = 44838
